# Fix FRAN CSV parsing: Use comma delimiter with quoted fields and proper US number format conversion
# Populate total_received, total_received_eur, automated_count, assigned_to_account,
# invoices_assigned, value_assigned, value_assigned_eur for rows 2-19 (previously all zero).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 817382.7100000001
$ws.Range("G2").Value = 817382.7100000001
$ws.Range("H2").Value = 43
$ws.Range("I2").Value = 54
$ws.Range("J2").Value = 163
$ws.Range("K2").Value = 430067.28
$ws.Range("L2").Value = 430067.28

$ws.Range("F3").Value = 887627.22
$ws.Range("G3").Value = 887627.22
$ws.Range("H3").Value = 98
$ws.Range("I3").Value = 133
$ws.Range("J3").Value = 136
$ws.Range("K3").Value = 325272.0600000001
$ws.Range("L3").Value = 325272.0600000001

$ws.Range("F4").Value = 8173894.139999999
$ws.Range("G4").Value = 9563456.143799998
$ws.Range("H4").Value = 31
$ws.Range("I4").Value = 48
$ws.Range("J4").Value = 104
$ws.Range("K4").Value = 691057.36
$ws.Range("L4").Value = 808537.1111999999

$ws.Range("F5").Value = 33735.12
$ws.Range("G5").Value = 33735.12
$ws.Range("H5").Value = 1
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 2
$ws.Range("K5").Value = 4096.1
$ws.Range("L5").Value = 4096.1

$ws.Range("F6").Value = 2520978.8
$ws.Range("G6").Value = 2520978.8
$ws.Range("H6").Value = 13
$ws.Range("I6").Value = 28
$ws.Range("J6").Value = 45
$ws.Range("K6").Value = 58858.24000000001
$ws.Range("L6").Value = 58858.24000000001

$ws.Range("F7").Value = 87712.25999999999
$ws.Range("G7").Value = 87712.25999999999
$ws.Range("H7").Value = 1
$ws.Range("I7").Value = 4
$ws.Range("J7").Value = 2
$ws.Range("K7").Value = 1420.71
$ws.Range("L7").Value = 1420.71

$ws.Range("F8").Value = 11601121.53
$ws.Range("G8").Value = 974494.2085199999
$ws.Range("H8").Value = 5
$ws.Range("I8").Value = 7
$ws.Range("J8").Value = 10
$ws.Range("K8").Value = 57023.46
$ws.Range("L8").Value = 4789.97064

$ws.Range("F9").Value = 649469.88
$ws.Range("G9").Value = 149378.0724
$ws.Range("H9").Value = 9
$ws.Range("I9").Value = 35
$ws.Range("J9").Value = 14
$ws.Range("K9").Value = 48199.26
$ws.Range("L9").Value = 11085.8298

$ws.Range("F10").Value = 5018728.44
$ws.Range("G10").Value = 5018728.44
$ws.Range("H10").Value = 49
$ws.Range("I10").Value = 82
$ws.Range("J10").Value = 208
$ws.Range("K10").Value = 1308416.98
$ws.Range("L10").Value = 1308416.98

$ws.Range("F11").Value = 230011.68
$ws.Range("G11").Value = 211610.7456
$ws.Range("H11").Value = 1
$ws.Range("I11").Value = 2
$ws.Range("J11").Value = 6
$ws.Range("K11").Value = 149673.18
$ws.Range("L11").Value = 137699.3256

$ws.Range("F12").Value = 1344304.89
$ws.Range("G12").Value = 1344304.89
$ws.Range("H12").Value = 1
$ws.Range("I12").Value = 45
$ws.Range("J12").Value = 1
$ws.Range("K12").Value = 663.23
$ws.Range("L12").Value = 663.23

$ws.Range("F13").Value = 353588.11
$ws.Range("G13").Value = 353588.11
$ws.Range("H13").Value = 7
$ws.Range("I13").Value = 359
$ws.Range("J13").Value = 7
$ws.Range("K13").Value = 6489.45
$ws.Range("L13").Value = 6489.45

$ws.Range("F14").Value = 5888894.439999999
$ws.Range("G14").Value = 5888894.439999999
$ws.Range("H14").Value = 30
$ws.Range("I14").Value = 63
$ws.Range("J14").Value = 66
$ws.Range("K14").Value = 112453.05
$ws.Range("L14").Value = 112453.05

$ws.Range("F15").Value = 204386.29
$ws.Range("G15").Value = 204386.29
$ws.Range("H15").Value = 7
$ws.Range("I15").Value = 11
$ws.Range("J15").Value = 12
$ws.Range("K15").Value = 11429.92
$ws.Range("L15").Value = 11429.92

$ws.Range("F16").Value = 2333841.4
$ws.Range("G16").Value = 2333841.4
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 1
$ws.Range("J16").Value = 0

$ws.Range("F17").Value = 129593.76
$ws.Range("G17").Value = 129593.76
$ws.Range("H17").Value = 2
$ws.Range("I17").Value = 3
$ws.Range("J17").Value = 2
$ws.Range("K17").Value = 918.2
$ws.Range("L17").Value = 918.2

$ws.Range("F18").Value = 7177062.160000001
$ws.Range("G18").Value = 7177062.160000001
$ws.Range("H18").Value = 68
$ws.Range("I18").Value = 108
$ws.Range("J18").Value = 247
$ws.Range("K18").Value = 474741.27
$ws.Range("L18").Value = 474741.27

$ws.Range("F19").Value = 2659261
$ws.Range("G19").Value = 103711.179
$ws.Range("H19").Value = 5
$ws.Range("I19").Value = 13
$ws.Range("J19").Value = 8
$ws.Range("K19").Value = 500526.8700000001
$ws.Range("L19").Value = 19520.54793
